$d = $word.ActiveDocument

$old = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$new = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Locate the credit line (it is spread across several differently
# formatted runs in the source document).
$rng = $d.Content
$found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Remove the whole run-fragmented credit text ...
    $rng.Delete()

    # ... and retype it as a single, plainly formatted run (matching
    # the target markup, which drops all the old per-run rPr/hyperlink
    # styling) with the astro-map year bumped from 2018 to 2022.
    $ins = $d.Range($start, $start)
    $ins.Font.Reset()
    $ins.InsertAfter($new)
}
